# Update column G ("K") values on Sheet1 to reflect the regenerated
# save_data using K (strike count) instead of Strike# for each row.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$kValues = @{
    2  = 8
    3  = 4
    4  = 7
    5  = 15
    6  = 6
    7  = 7
    8  = 4
    9  = 7
    10 = 7
    11 = 5
    12 = 2
    13 = 5
    14 = 3
    15 = 2
    16 = 5
    17 = 3
    18 = 1
    19 = 5
    20 = 3
    21 = 6
    22 = 4
    23 = 3
    24 = 0
    25 = 2
    26 = 2
    27 = 2
    28 = 0
    29 = 2
    30 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
